$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$value)
    if ($value -match '^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$') {
        $cell.NumberFormat = "@"
        $cell.Value = $value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $value
    }
}

Set-TextValue $ws.Cells.Item(2, 4) "30.191.49"
$ws.Cells.Item(2, 5).Value = "  +0.04%  "

Set-TextValue $ws.Cells.Item(3, 4) "1.868.31"
$ws.Cells.Item(3, 5).Value = "  +0.79%  "

Set-TextValue $ws.Cells.Item(4, 4) "1.001"
$ws.Cells.Item(4, 5).Value = "  +0.09%  "

Set-TextValue $ws.Cells.Item(5, 4) "234.69"
$ws.Cells.Item(5, 5).Value = "  -0.27%  "

Set-TextValue $ws.Cells.Item(6, 4) "1.001"
$ws.Cells.Item(6, 5).Value = "  +0.09%  "

Set-TextValue $ws.Cells.Item(7, 4) "0.4697"
$ws.Cells.Item(7, 5).Value = "  +0.12%  "

Set-TextValue $ws.Cells.Item(8, 4) "0.2847"
$ws.Cells.Item(8, 5).Value = "  -1.42%  "

Set-TextValue $ws.Cells.Item(9, 4) "41.49"
$ws.Cells.Item(9, 5).Value = "  -2.79%  "

Set-TextValue $ws.Cells.Item(10, 4) "0.06561"
$ws.Cells.Item(10, 5).Value = "  +0.10%  "

Set-TextValue $ws.Cells.Item(11, 4) "21.21"
$ws.Cells.Item(11, 5).Value = "  -2.79%  "

Set-TextValue $ws.Cells.Item(12, 4) "0.07783"
$ws.Cells.Item(12, 5).Value = "  -2.37%  "

Set-TextValue $ws.Cells.Item(13, 4) "96.18"
$ws.Cells.Item(13, 5).Value = "  -1.30%  "

Set-TextValue $ws.Cells.Item(14, 4) "1.869.11"
$ws.Cells.Item(14, 5).Value = "  +0.86%  "

Set-TextValue $ws.Cells.Item(15, 4) "0.6900"
$ws.Cells.Item(15, 5).Value = "  +2.08%  "

Set-TextValue $ws.Cells.Item(16, 4) "5.098"
$ws.Cells.Item(16, 5).Value = "  -0.14%  "

Set-TextValue $ws.Cells.Item(17, 4) "266.25"
$ws.Cells.Item(17, 5).Value = "  -0.87%  "

Set-TextValue $ws.Cells.Item(18, 4) "30.193.21"
$ws.Cells.Item(18, 5).Value = "  +0.15%  "

Set-TextValue $ws.Cells.Item(19, 4) "13.67"
$ws.Cells.Item(19, 5).Value = "  +0.39%  "

Set-TextValue $ws.Cells.Item(20, 4) "0.000007726"
$ws.Cells.Item(20, 5).Value = "  +0.25%  "

$ws.Cells.Item(21, 5).Value = "  -0.03%  "

Set-TextValue $ws.Cells.Item(22, 4) "2.123.47"
$ws.Cells.Item(22, 5).Value = "  +1.38%  "

Set-TextValue $ws.Cells.Item(23, 4) "1.001"
$ws.Cells.Item(23, 5).Value = "  +0.04%  "

Set-TextValue $ws.Cells.Item(24, 4) "5.228"
$ws.Cells.Item(24, 5).Value = "  +0.71%  "

Set-TextValue $ws.Cells.Item(25, 4) "6.158"
$ws.Cells.Item(25, 5).Value = "  +0.37%  "

Set-TextValue $ws.Cells.Item(26, 4) "9.467"
$ws.Cells.Item(26, 5).Value = "  +3.36%  "

Set-TextValue $ws.Cells.Item(27, 4) "166.07"
$ws.Cells.Item(27, 5).Value = "  -0.64%  "

Set-TextValue $ws.Cells.Item(28, 4) "18.72"
$ws.Cells.Item(28, 5).Value = "  -0.77%  "

$ws.Cells.Item(29, 5).Value = "  +0.36%  "

Set-TextValue $ws.Cells.Item(30, 4) "1.372"
$ws.Cells.Item(30, 5).Value = "  -0.38%  "

Set-TextValue $ws.Cells.Item(31, 4) "0.09932"
$ws.Cells.Item(31, 5).Value = "  +0.78%  "

Set-TextValue $ws.Cells.Item(32, 4) "4.356"
$ws.Cells.Item(32, 5).Value = "  +1.62%  "

Set-TextValue $ws.Cells.Item(33, 4) "1.459"
$ws.Cells.Item(33, 5).Value = "  -0.21%  "

Set-TextValue $ws.Cells.Item(34, 4) "4.042"
$ws.Cells.Item(34, 5).Value = "  +1.15%  "

Set-TextValue $ws.Cells.Item(35, 4) "0.04728"
$ws.Cells.Item(35, 5).Value = "  +0.69%  "

Set-TextValue $ws.Cells.Item(36, 4) "1.129"
$ws.Cells.Item(36, 5).Value = "  +1.02%  "

Set-TextValue $ws.Cells.Item(37, 4) "0.7002"
$ws.Cells.Item(37, 5).Value = "  +0.25%  "

Set-TextValue $ws.Cells.Item(38, 4) "2.715"
$ws.Cells.Item(38, 5).Value = "  +0.35%  "

Set-TextValue $ws.Cells.Item(39, 4) "0.01862"
$ws.Cells.Item(39, 5).Value = "  -0.36%  "

Set-TextValue $ws.Cells.Item(40, 4) "2.774"
$ws.Cells.Item(40, 5).Value = "  +6.65%  "

Set-TextValue $ws.Cells.Item(41, 4) "6.239"
$ws.Cells.Item(41, 5).Value = "  -1.28%  "

Set-TextValue $ws.Cells.Item(42, 4) "72.65"
$ws.Cells.Item(42, 5).Value = "  -0.76%  "

Set-TextValue $ws.Cells.Item(43, 4) "1.936"
$ws.Cells.Item(43, 5).Value = "  +0.17%  "

$ws.Cells.Item(44, 5).Value = "  +0.14%  "

$ws.Cells.Item(45, 5).Value = "  +0.38%  "

Set-TextValue $ws.Cells.Item(46, 4) "0.8338"
$ws.Cells.Item(46, 5).Value = "  -0.60%  "

Set-TextValue $ws.Cells.Item(47, 4) "102.68"
$ws.Cells.Item(47, 5).Value = "  -0.41%  "

Set-TextValue $ws.Cells.Item(48, 4) "972.10"
$ws.Cells.Item(48, 5).Value = "  +3.84%  "

Set-TextValue $ws.Cells.Item(49, 4) "7.064"
$ws.Cells.Item(49, 5).Value = "  +0.66%  "

Set-TextValue $ws.Cells.Item(50, 4) "9.166"
$ws.Cells.Item(50, 5).Value = "  +0.39%  "

Set-TextValue $ws.Cells.Item(51, 4) "34.48"
$ws.Cells.Item(51, 5).Value = "  +1.92%  "
